$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers so the "old"/"new" suffixes become the actual
#    format-version identifiers used for this comparison (FV2210 / FV2304).
# ---------------------------------------------------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ([string]$cell.Value2) -replace "_old$", "_FV2210"
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ([string]$cell.Value2) -replace "_new$", "_FV2304"
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an actual Excel Table (ListObject) so the
#    headers also drive the table's column definitions/autofilter.
#    The header row already carries bespoke (bold/shaded/bordered) direct
#    formatting; converting straight into a table would make Excel snapshot
#    that look into a new header-row dxf, so we stash it, wipe it, build the
#    table against now-plain header cells, then paint the very same look
#    back on afterwards (this matches how the reference workbook looks -
#    no extra dxf, no style bloat).
# ---------------------------------------------------------------------------
$headerRng = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")

$headerRng.Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$headerRng.ClearFormats()

$tableRng = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $tableRng, $null, 1)
$tbl.Name = "Table1"

$scratch.Copy()
$headerRng.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# fully remove the scratch row again (not just clear it) so it does not
# linger as a blank row and inflate the sheet dimension
$scratch.EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3) Freeze the header row so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
